$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.966392159461975
$ws.Range("B1").Value = 3.504112005233765
$ws.Range("C1").Value = 2.3186936378479
$ws.Range("D1").Value = 2.024323463439941
$ws.Range("E1").Value = 1.86319625377655
